# Asset List.xlsx edit: progress through reaper and fmod
#
# Replace "reaper production" (shared string) with a new "FMOD" shared
# string on the rows that have moved from Reaper to FMOD, and update the
# active cell selection on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmodRows = @(5, 11, 12, 17, 18, 19, 20, 23, 24, 30)
foreach ($row in $fmodRows) {
    $ws.Range("C$row").Value = "FMOD"
}

# Update the selected cell on the sheet to reflect where the author left off.
$ws.Range("C21").Select()
